# Update the "Förändrad" (Changed) date column from 2023-09-16 (45185)
# to 2023-10-05 (45204) for the data rows in the "Avverkningsanmälningar" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

$ws.Range("C2:C15").Value = 45204
